# Scheduled runner refresh: update market-price derived columns (H-N)
# across the Leve profit tables on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 342.13925
$ws.Range("J17").Value = 342.13925
$ws.Range("L17").Value = 1026.41775
$ws.Range("N17").Value = -1362.41775
$ws.Range("H64").Value = 3071.4285
$ws.Range("I64").Value = 2825
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 2825
$ws.Range("L64").Value = 3400
$ws.Range("M64").Value = -2577
$ws.Range("N64").Value = -3896
$ws.Range("H67").Value = 3071.4285
$ws.Range("I67").Value = 2825
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 2825
$ws.Range("L67").Value = 3400
$ws.Range("M67").Value = -1967
$ws.Range("N67").Value = -5116
$ws.Range("H112").Value = 956.96826
$ws.Range("I112").Value = 425
$ws.Range("J112").Value = 993.0339
$ws.Range("K112").Value = 1275
$ws.Range("L112").Value = 2979.1017
$ws.Range("M112").Value = -167
$ws.Range("N112").Value = -5195.1017
$ws.Range("H138").Value = 2301.5173
$ws.Range("I138").Value = 1402.6538
$ws.Range("J138").Value = 3031.8438
$ws.Range("K138").Value = 4207.9614
$ws.Range("L138").Value = 9095.5314
$ws.Range("M138").Value = 932.0385999999999
$ws.Range("N138").Value = -19375.5314
$ws.Range("H141").Value = 5619.2
$ws.Range("I141").Value = 2698.5
$ws.Range("J141").Value = 10000.25
$ws.Range("K141").Value = 8095.5
$ws.Range("L141").Value = 30000.75
$ws.Range("M141").Value = -2915.5
$ws.Range("N141").Value = -40360.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1147.7941
$ws.Range("I74").Value = 1067.5
$ws.Range("K74").Value = 1067.5
$ws.Range("M74").Value = -193.5
$ws.Range("H77").Value = 1147.7941
$ws.Range("I77").Value = 1067.5
$ws.Range("K77").Value = 5337.5
$ws.Range("M77").Value = -969.5
$ws.Range("H110").Value = 1415.1034
$ws.Range("I110").Value = 973
$ws.Range("J110").Value = 2804.5715
$ws.Range("K110").Value = 973
$ws.Range("L110").Value = 2804.5715
$ws.Range("M110").Value = 1072
$ws.Range("N110").Value = -6894.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 553.1389
$ws.Range("I94").Value = 518.3214
$ws.Range("J94").Value = 675
$ws.Range("K94").Value = 518.3214
$ws.Range("L94").Value = 675
$ws.Range("M94").Value = -67.32140000000004
$ws.Range("N94").Value = -1577
$ws.Range("H99").Value = 1523.9412
$ws.Range("I99").Value = 1138.081
$ws.Range("K99").Value = 1138.081
$ws.Range("M99").Value = 359.9190000000001
$ws.Range("H134").Value = 2181816.2
$ws.Range("I134").Value = 1173.1464
$ws.Range("J134").Value = 11122454
$ws.Range("K134").Value = 3519.4392
$ws.Range("L134").Value = 33367362
$ws.Range("M134").Value = -984.4392000000003
$ws.Range("N134").Value = -33372432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1326.3636
$ws.Range("I16").Value = 851.4
$ws.Range("K16").Value = 851.4
$ws.Range("M16").Value = -564.4
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 90912776
$ws.Range("I99").Value = 250002500
$ws.Range("J99").Value = 4359.143
$ws.Range("K99").Value = 250002500
$ws.Range("L99").Value = 4359.143
$ws.Range("M99").Value = -250001002
$ws.Range("N99").Value = -7355.143
$ws.Range("H113").Value = 1326.3636
$ws.Range("I113").Value = 851.4
$ws.Range("K113").Value = 851.4
$ws.Range("M113").Value = 1318.6
$ws.Range("H122").Value = 22728028
$ws.Range("I122").Value = 41667164
$ws.Range("K122").Value = 125001492
$ws.Range("M122").Value = -124999042
$ws.Range("H126").Value = 90912776
$ws.Range("I126").Value = 250002500
$ws.Range("J126").Value = 4359.143
$ws.Range("K126").Value = 750007500
$ws.Range("L126").Value = 13077.429
$ws.Range("M126").Value = -750005030
$ws.Range("N126").Value = -18017.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 14186927
$ws.Range("I5").Value = 19047908
$ws.Range("J5").Value = 9062.5
$ws.Range("K5").Value = 57143724
$ws.Range("L5").Value = 27187.5
$ws.Range("M5").Value = -57143612
$ws.Range("N5").Value = -27411.5
$ws.Range("H114").Value = 709.5
$ws.Range("I114").Value = 337.4
$ws.Range("J114").Value = 790.3913
$ws.Range("K114").Value = 1012.2
$ws.Range("L114").Value = 2371.1739
$ws.Range("M114").Value = 2241.8
$ws.Range("N114").Value = -8879.1739
$ws.Range("H122").Value = 12259502
$ws.Range("I122").Value = 78125336
$ws.Range("J122").Value = 5393.6978
$ws.Range("K122").Value = 703128024
$ws.Range("L122").Value = 48543.2802
$ws.Range("M122").Value = -703125574
$ws.Range("N122").Value = -53443.2802
$ws.Range("H134").Value = 13890913
$ws.Range("I134").Value = 17858318
$ws.Range("J134").Value = 4998.75
$ws.Range("K134").Value = 53574954
$ws.Range("L134").Value = 14996.25
$ws.Range("M134").Value = -53569884
$ws.Range("N134").Value = -25136.25
$ws.Range("H135").Value = 14186927
$ws.Range("I135").Value = 19047908
$ws.Range("J135").Value = 9062.5
$ws.Range("K135").Value = 171431172
$ws.Range("L135").Value = 81562.5
$ws.Range("M135").Value = -171428637
$ws.Range("N135").Value = -86632.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1453
$ws.Range("I61").Value = 1408.1818
$ws.Range("J61").Value = 1576.25
$ws.Range("K61").Value = 1408.1818
$ws.Range("L61").Value = 1576.25
$ws.Range("M61").Value = -1206.1818
$ws.Range("N61").Value = -1980.25
$ws.Range("H82").Value = 1056.96
$ws.Range("I82").Value = 874.4286
$ws.Range("J82").Value = 1289.2727
$ws.Range("K82").Value = 874.4286
$ws.Range("L82").Value = 1289.2727
$ws.Range("M82").Value = -513.4286
$ws.Range("N82").Value = -2011.2727
$ws.Range("H85").Value = 1056.96
$ws.Range("I85").Value = 874.4286
$ws.Range("J85").Value = 1289.2727
$ws.Range("K85").Value = 874.4286
$ws.Range("L85").Value = 1289.2727
$ws.Range("M85").Value = 373.5714
$ws.Range("N85").Value = -3785.2727
$ws.Range("H113").Value = 1453
$ws.Range("I113").Value = 1408.1818
$ws.Range("J113").Value = 1576.25
$ws.Range("K113").Value = 1408.1818
$ws.Range("L113").Value = 1576.25
$ws.Range("M113").Value = 761.8181999999999
$ws.Range("N113").Value = -5916.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28811.947
$ws.Range("I122").Value = 43202.168
$ws.Range("J122").Value = 4143
$ws.Range("K122").Value = 129606.504
$ws.Range("L122").Value = 12429
$ws.Range("M122").Value = -127156.504
$ws.Range("N122").Value = -17329
$ws.Range("H126").Value = 1154.5
$ws.Range("I126").Value = 713
$ws.Range("J126").Value = 1596
$ws.Range("K126").Value = 2139
$ws.Range("L126").Value = 4788
$ws.Range("M126").Value = 331
$ws.Range("N126").Value = -9728
$ws.Range("H132").Value = 7264098
$ws.Range("I132").Value = 23731.088
$ws.Range("K132").Value = 71193.264
$ws.Range("M132").Value = -68663.264
